$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("G5").Value = 1.8
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 2.5
$ws.Range("K5").Value = 2.05
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 2.75
$ws.Range("Q5").Value = 2.2
$ws.Range("R5").Value = 1.65
$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63
$ws.Range("U5").Value = 2.1
$ws.Range("V5").Value = 1.67
$ws.Range("W5").Value = 6
$ws.Range("X5").Value = 7.5
$ws.Range("AA5").Value = 17
$ws.Range("AC5").Value = 8.5
$ws.Range("AF5").Value = 67
$ws.Range("AG5").Value = 9.5
$ws.Range("AH5").Value = 21
$ws.Range("AM5").Value = 900
$ws.Range("AO5").Value = 10
$ws.Range("AP5").Value = 23
$ws.Range("AQ5").Value = 34
$ws.Range("AR5").Value = 51
$ws.Range("AT5").Value = 2.63
$ws.Range("AY5").Value = 34
$ws.Range("BB5").Value = 301

# Row 7
$ws.Range("G7").Value = 1.37
$ws.Range("H7").Value = 4.55
$ws.Range("I7").Value = 7.2
$ws.Range("K7").Value = 2.45
$ws.Range("L7").Value = 6.3
$ws.Range("P7").Value = 4.25
$ws.Range("Q7").Value = 1.53
$ws.Range("R7").Value = 2.18
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 1.83
$ws.Range("W7").Value = 7.9
$ws.Range("X7").Value = 7
$ws.Range("AA7").Value = 10.75
$ws.Range("AB7").Value = 24
$ws.Range("AC7").Value = 14.5
$ws.Range("AD7").Value = 9.25
$ws.Range("AK7").Value = 75
$ws.Range("AN7").Value = 3.25
$ws.Range("AO7").Value = 6.1
$ws.Range("AT7").Value = 3.25
$ws.Range("AU7").Value = 7.8
$ws.Range("AY7").Value = 35

# Row 10
$ws.Range("G10").Value = 4.33
$ws.Range("I10").Value = 1.73
$ws.Range("J10").Value = 5.5
$ws.Range("L10").Value = 2.4
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 8
$ws.Range("Q10").Value = 2.2
$ws.Range("R10").Value = 1.65
$ws.Range("U10").Value = 2.1
$ws.Range("V10").Value = 1.67
$ws.Range("W10").Value = 11
$ws.Range("X10").Value = 23
$ws.Range("Z10").Value = 51
$ws.Range("AA10").Value = 41
$ws.Range("AB10").Value = 51
$ws.Range("AC10").Value = 8
$ws.Range("AE10").Value = 19
$ws.Range("AG10").Value = 6
$ws.Range("AH10").Value = 7.5
$ws.Range("AJ10").Value = 13
$ws.Range("AN10").Value = 6.5
$ws.Range("AO10").Value = 29
$ws.Range("AP10").Value = 41
$ws.Range("AQ10").Value = 101
$ws.Range("AS10").Value = 351
$ws.Range("AU10").Value = 9
$ws.Range("AW10").Value = 3.6
$ws.Range("AX10").Value = 9.5

# Row 12
$ws.Range("H12").Value = 3.25
$ws.Range("I12").Value = 2.88
$ws.Range("K12").Value = 2.1
$ws.Range("M12").Value = 1.05
$ws.Range("N12").Value = 11
$ws.Range("O12").Value = 1.29
$ws.Range("P12").Value = 3.5
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.85
$ws.Range("S12").Value = 1.4
$ws.Range("T12").Value = 2.75
$ws.Range("U12").Value = 1.8
$ws.Range("V12").Value = 1.91
$ws.Range("W12").Value = 8
$ws.Range("Z12").Value = 21
$ws.Range("AA12").Value = 19
$ws.Range("AB12").Value = 29
$ws.Range("AC12").Value = 9.5
$ws.Range("AG12").Value = 9
$ws.Range("AK12").Value = 23
$ws.Range("AM12").Value = 251
$ws.Range("AS12").Value = 151
$ws.Range("AT12").Value = 2.75
$ws.Range("AY12").Value = 26

# Row 13
$ws.Range("G13").Value = 1.91
$ws.Range("I13").Value = 3.8
$ws.Range("J13").Value = 2.63
$ws.Range("L13").Value = 4.33
$ws.Range("W13").Value = 6.5
$ws.Range("X13").Value = 9
$ws.Range("Z13").Value = 17
$ws.Range("AA13").Value = 17
$ws.Range("AD13").Value = 6.5
$ws.Range("AG13").Value = 10
$ws.Range("AH13").Value = 19
$ws.Range("AK13").Value = 34
$ws.Range("AO13").Value = 11
$ws.Range("AR13").Value = 51
$ws.Range("BB13").Value = 251

# Row 20
$ws.Range("N20").Value = 8

# Row 21
$ws.Range("Q21").Value = 2.5
$ws.Range("R21").Value = 1.5

# Row 22
$ws.Range("G22").Value = 6.5
$ws.Range("H22").Value = 3.5
$ws.Range("I22").Value = 1.57
$ws.Range("J22").Value = 6.5
$ws.Range("AD22").Value = 7
$ws.Range("AF22").Value = 67
$ws.Range("AJ22").Value = 11
$ws.Range("AO22").Value = 34
$ws.Range("AR22").Value = 151
$ws.Range("AW22").Value = 3.4

# Row 23
$ws.Range("L23").Value = 2.87

# Row 27
$ws.Range("H27").Value = 3
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 2.67
$ws.Range("K27").Value = 1.91
$ws.Range("L27").Value = 4.5
$ws.Range("Q27").Value = 2.35
$ws.Range("S27").Value = 1.52
$ws.Range("T27").Value = 2.22
$ws.Range("U27").Value = 2.05
$ws.Range("W27").Value = 5.4
$ws.Range("X27").Value = 8
$ws.Range("Y27").Value = 9
$ws.Range("AA27").Value = 20
$ws.Range("AB27").Value = 40
$ws.Range("AC27").Value = 6.5
$ws.Range("AG27").Value = 9
$ws.Range("AH27").Value = 21
$ws.Range("AI27").Value = 14
$ws.Range("AJ27").Value = 65
$ws.Range("AK27").Value = 45
$ws.Range("AL27").Value = 60
$ws.Range("AN27").Value = 3.65
$ws.Range("AO27").Value = 10.75
$ws.Range("AP27").Value = 24
$ws.Range("AQ27").Value = 45
$ws.Range("AR27").Value = 100
$ws.Range("AS27").Value = 400
$ws.Range("AT27").Value = 2.2
$ws.Range("AW27").Value = 5.5
$ws.Range("AX27").Value = 24
$ws.Range("AY27").Value = 32
$ws.Range("AZ27").Value = 150
$ws.Range("BA27").Value = 200
$ws.Range("BB27").Value = 500
